# Auto-generated edit script: update F-column attendance/sales counts
# per commit 456a3b4 (gh-pages data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 506  # was 505
$ws.Range("F3").Value = 1582  # was 1579
$ws.Range("F4").Value = 826  # was 823
$ws.Range("F6").Value = 66  # was 64
$ws.Range("F7").Value = 1139  # was 1134
$ws.Range("F8").Value = 743  # was 738
$ws.Range("F9").Value = 790  # was 783
$ws.Range("F10").Value = 1427  # was 1418
$ws.Range("F11").Value = 285  # was 280
$ws.Range("F12").Value = 1030  # was 1028
$ws.Range("F13").Value = 31  # was 30
$ws.Range("F15").Value = 191  # was 189
$ws.Range("F16").Value = 49  # was 48
$ws.Range("F17").Value = 461  # was 454
$ws.Range("F18").Value = 12  # was 4
$ws.Range("F19").Value = 23  # was 18
$ws.Range("F20").Value = 2  # was 0
$ws.Range("F22").Value = 294  # was 293
$ws.Range("F23").Value = 550  # was 547
$ws.Range("F24").Value = 564  # was 560
$ws.Range("F26").Value = 242  # was 240
$ws.Range("F27").Value = 176  # was 174
$ws.Range("F28").Value = 368  # was 367
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 226  # was 222
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 506  # was 505
$ws.Range("F3").Value = 226  # was 222
$ws.Range("F4").Value = 1582  # was 1579
$ws.Range("F6").Value = 826  # was 823
$ws.Range("F9").Value = 66  # was 64
$ws.Range("F10").Value = 1139  # was 1134
$ws.Range("F11").Value = 743  # was 738
$ws.Range("F12").Value = 790  # was 783
$ws.Range("F13").Value = 1427  # was 1418
$ws.Range("F14").Value = 285  # was 281
$ws.Range("F15").Value = 1030  # was 1028
$ws.Range("F16").Value = 31  # was 30
$ws.Range("F18").Value = 191  # was 189
$ws.Range("F19").Value = 49  # was 48
$ws.Range("F20").Value = 461  # was 454
$ws.Range("F21").Value = 12  # was 4
$ws.Range("F22").Value = 23  # was 18
$ws.Range("F24").Value = 2  # was 0
$ws.Range("F27").Value = 294  # was 293
$ws.Range("F31").Value = 550  # was 547
$ws.Range("F32").Value = 564  # was 560
$ws.Range("F34").Value = 242  # was 240
$ws.Range("F36").Value = 176  # was 174
$ws.Range("F41").Value = 368  # was 367
